$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reshuffles the data rows (2-6) of the sheet: each row's content
# (columns A, B, E, F, G, H, Q, R and AC) moves to a different row while the
# rest of the row (dates, observer names, etc.) stays put. Apply the new
# values for columns A, B, E, F, G, H, Q, R, AC directly.

$ws.Range("A2").Value = 111638281
$ws.Range("B2").Value = 89423
$ws.Range("E2").Value = 5432
$ws.Range("F2").Value = "Granticka"
$ws.Range("G2").Value = "Porodaedalea chrysoloma"
$ws.Range("H2").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q2").Value = 588278.4784540196
$ws.Range("R2").Value = 7033316.716959814
$ws.Range("AC2").Value = ""

$ws.Range("A3").Value = 111638278
$ws.Range("B3").Value = 56414
$ws.Range("E3").Value = 100049
$ws.Range("F3").Value = "Spillkråka"
$ws.Range("G3").Value = "Dryocopus martius"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("Q3").Value = 588319.7233830886
$ws.Range("R3").Value = 7033285.591169797
$ws.Range("AC3").Value = "Bohål i gammal grov tall."

$ws.Range("A4").Value = 111638283
$ws.Range("B4").Value = 77515
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 588212.5872669512
$ws.Range("R4").Value = 7033297.86989607
$ws.Range("AC4").Value = ""

$ws.Range("A5").Value = 111638277
$ws.Range("B5").Value = 77267
$ws.Range("E5").Value = 6446
$ws.Range("F5").Value = "Kolflarnlav"
$ws.Range("G5").Value = "Carbonicola anthracophila"
$ws.Range("H5").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q5").Value = 588323.0816159723
$ws.Range("R5").Value = 7033261.485134664
$ws.Range("AC5").Value = "På kolad tallstubbe med yxhugg."

$ws.Range("A6").Value = 111638282
$ws.Range("B6").Value = 89405
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = "Ullticka"
$ws.Range("G6").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 588265.5989708689
$ws.Range("R6").Value = 7033311.880202802
$ws.Range("AC6").Value = ""
